$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (ECs as target cluster / duplicate MuSCs-sourced rows removed in new data)
$ws.Range("A8:T10").Delete()

$row2 = @("ECs", "Cxcl17", "Gpr35", "FAPs", 1, 0.3333333333333333, 0.09427400000000001, 0.282822, 0.5730152684128631, 0.5730152684128631, 3, 1, 0.6379899999999999, 1.91397, 0.8729739184691609, 0.8729739184691611, 0.06014586926, 0.54131282334, 0.5002273842090351, 0.5002273842090352)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @("ECs", "Cxcl17", "Gpr35", "MuSCs", 1, 0.3333333333333333, 0.09427400000000001, 0.282822, 0.5730152684128631, 0.5730152684128631, 1, 0.3333333333333333, 0.09283366666666666, 0.278501, 0.1270260815308389, 0.127026081530839, 0.008751801091333333, 0.078766209822, 0.07278788420382791, 0.07278788420382792)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @("FAPs", "Cxcl17", "Gpr35", "FAPs", 1, 0.3333333333333333, 0.023327, 0.069981, 0.1417859342582987, 0.1417859342582987, 3, 1, 0.6379899999999999, 1.91397, 0.8729739184691609, 0.8729739184691611, 0.01488239273, 0.13394153457, 0.1237754226132779, 0.1237754226132779)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @("FAPs", "Cxcl17", "Gpr35", "MuSCs", 1, 0.3333333333333333, 0.023327, 0.069981, 0.1417859342582987, 0.1417859342582987, 1, 0.3333333333333333, 0.09283366666666666, 0.278501, 0.1270260815308389, 0.127026081530839, 0.002165530942333333, 0.019489778481, 0.01801051164502083, 0.01801051164502083)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

$row6 = @("MuSCs", "Cxcl17", "Gpr35", "FAPs", 1, 0.3333333333333333, 0.04692166666666667, 0.140765, 0.2851987973288382, 0.2851987973288382, 3, 1, 0.6379899999999999, 1.91397, 0.8729739184691609, 0.8729739184691611, 0.02993555411666667, 0.26941998705, 0.2489711116468479, 0.248971111646848)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6[$i]
}

$row7 = @("MuSCs", "Cxcl17", "Gpr35", "MuSCs", 1, 0.3333333333333333, 0.04692166666666667, 0.140765, 0.2851987973288382, 0.2851987973288382, 1, 0.3333333333333333, 0.09283366666666666, 0.278501, 0.1270260815308389, 0.127026081530839, 0.004355910362777778, 0.039203193265, 0.03622768568199021, 0.03622768568199022)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $row7[$i]
}
